$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1149.4615
$ws.Range("I12").Value = 2256.8333
$ws.Range("J12").Value = 200.28572
$ws.Range("K12").Value = 2256.8333
$ws.Range("L12").Value = 200.28572
$ws.Range("M12").Value = -2086.8333
$ws.Range("N12").Value = -540.28572
$ws.Range("H40").Value = 1372.7
$ws.Range("I40").Value = 787.5
$ws.Range("J40").Value = 1762.8334
$ws.Range("K40").Value = 787.5
$ws.Range("L40").Value = 1762.8334
$ws.Range("M40").Value = -612.5
$ws.Range("N40").Value = -2112.8334
$ws.Range("H64").Value = 2965.6
$ws.Range("I64").Value = 2815.1428
$ws.Range("K64").Value = 2815.1428
$ws.Range("M64").Value = -2567.1428
$ws.Range("H67").Value = 2965.6
$ws.Range("I67").Value = 2815.1428
$ws.Range("K67").Value = 2815.1428
$ws.Range("M67").Value = -1957.1428
$ws.Range("H76").Value = 3706683.8
$ws.Range("I76").Value = 3100
$ws.Range("J76").Value = 5053441.5
$ws.Range("K76").Value = 3100
$ws.Range("L76").Value = 5053441.5
$ws.Range("M76").Value = -2785
$ws.Range("N76").Value = -5054071.5
$ws.Range("H79").Value = 3706683.8
$ws.Range("I79").Value = 3100
$ws.Range("J79").Value = 5053441.5
$ws.Range("K79").Value = 3100
$ws.Range("L79").Value = 5053441.5
$ws.Range("M79").Value = -2008
$ws.Range("N79").Value = -5055625.5
$ws.Range("H129").Value = 769.7636
$ws.Range("J129").Value = 800.7843
$ws.Range("L129").Value = 2402.3529
$ws.Range("N129").Value = -12402.3529
$ws.Range("H132").Value = 3279.606
$ws.Range("I132").Value = 3420.4075
$ws.Range("K132").Value = 10261.2225
$ws.Range("M132").Value = -7731.2225
$ws.Range("H137").Value = 83108.86
$ws.Range("I137").Value = 96253.12
$ws.Range("K137").Value = 288759.36
$ws.Range("M137").Value = -286209.36
$ws.Range("H141").Value = 2044.238
$ws.Range("I141").Value = 1683.1765
$ws.Range("J141").Value = 3578.75
$ws.Range("K141").Value = 5049.529500000001
$ws.Range("L141").Value = 10736.25
$ws.Range("M141").Value = 130.4704999999994
$ws.Range("N141").Value = -21096.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7217.1646
$ws.Range("I32").Value = 5819.8237
$ws.Range("J32").Value = 12806.529
$ws.Range("K32").Value = 5819.8237
$ws.Range("L32").Value = 12806.529
$ws.Range("M32").Value = -5532.8237
$ws.Range("N32").Value = -13380.529
$ws.Range("H61").Value = 2899.2856
$ws.Range("I61").Value = 3432.88
$ws.Range("J61").Value = 2114.5881
$ws.Range("K61").Value = 3432.88
$ws.Range("L61").Value = 2114.5881
$ws.Range("M61").Value = -3220.88
$ws.Range("N61").Value = -2538.5881
$ws.Range("H63").Value = 3908702.5
$ws.Range("I63").Value = 2802.7144
$ws.Range("K63").Value = 2802.7144
$ws.Range("M63").Value = -2116.7144
$ws.Range("H66").Value = 3908702.5
$ws.Range("I66").Value = 2802.7144
$ws.Range("K66").Value = 14013.572
$ws.Range("M66").Value = -10581.572
$ws.Range("H88").Value = 334618
$ws.Range("I88").Value = 1900
$ws.Range("K88").Value = 1900
$ws.Range("M88").Value = -1494
$ws.Range("H91").Value = 334618
$ws.Range("I91").Value = 1900
$ws.Range("K91").Value = 1900
$ws.Range("M91").Value = -496
$ws.Range("H122").Value = 1755.1892
$ws.Range("I122").Value = 1651.7667
$ws.Range("K122").Value = 4955.300099999999
$ws.Range("M122").Value = -2505.300099999999
$ws.Range("H132").Value = 16015.838
$ws.Range("I132").Value = 2388.3914
$ws.Range("J132").Value = 38403.785
$ws.Range("K132").Value = 7165.174199999999
$ws.Range("L132").Value = 115211.355
$ws.Range("M132").Value = -4635.174199999999
$ws.Range("N132").Value = -120271.355
$ws.Range("H136").Value = 2899.2856
$ws.Range("I136").Value = 3432.88
$ws.Range("J136").Value = 2114.5881
$ws.Range("K136").Value = 10298.64
$ws.Range("L136").Value = 6343.7643
$ws.Range("M136").Value = -7748.639999999999
$ws.Range("N136").Value = -11443.7643

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 911002.75
$ws.Range("I105").Value = 1502.6086
$ws.Range("K105").Value = 1502.6086
$ws.Range("M105").Value = 244.3914
$ws.Range("H134").Value = 4169.2666
$ws.Range("I134").Value = 4352.2856
$ws.Range("J134").Value = 1607
$ws.Range("K134").Value = 13056.8568
$ws.Range("L134").Value = 4821
$ws.Range("M134").Value = -10521.8568
$ws.Range("N134").Value = -9891

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4841.927
$ws.Range("I31").Value = 2850.8462
$ws.Range("J31").Value = 5766.357
$ws.Range("K31").Value = 2850.8462
$ws.Range("L31").Value = 5766.357
$ws.Range("M31").Value = -2555.8462
$ws.Range("N31").Value = -6356.357
$ws.Range("H34").Value = 4841.927
$ws.Range("I34").Value = 2850.8462
$ws.Range("J34").Value = 5766.357
$ws.Range("K34").Value = 2850.8462
$ws.Range("L34").Value = 5766.357
$ws.Range("M34").Value = -2648.8462
$ws.Range("N34").Value = -6170.357
$ws.Range("H58").Value = 52576.7
$ws.Range("J58").Value = 64256.875
$ws.Range("L58").Value = 64256.875
$ws.Range("N58").Value = -64662.875
$ws.Range("H62").Value = 6668.6665
$ws.Range("J62").Value = 5003
$ws.Range("L62").Value = 5003
$ws.Range("N62").Value = -6251
$ws.Range("H65").Value = 6668.6665
$ws.Range("J65").Value = 5003
$ws.Range("L65").Value = 25015
$ws.Range("N65").Value = -31255
$ws.Range("H99").Value = 26319708
$ws.Range("I99").Value = 3054
$ws.Range("K99").Value = 3054
$ws.Range("M99").Value = -1556
$ws.Range("H126").Value = 26319708
$ws.Range("I126").Value = 3054
$ws.Range("K126").Value = 9162
$ws.Range("M126").Value = -6692
$ws.Range("H132").Value = 4214
$ws.Range("I132").Value = 3348.4443
$ws.Range("J132").Value = 5187.75
$ws.Range("K132").Value = 10045.3329
$ws.Range("L132").Value = 15563.25
$ws.Range("M132").Value = -7515.332900000001
$ws.Range("N132").Value = -20623.25
$ws.Range("H136").Value = 52576.7
$ws.Range("J136").Value = 64256.875
$ws.Range("L136").Value = 192770.625
$ws.Range("N136").Value = -197870.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1974.75
$ws.Range("I69").Value = 1799.6666
$ws.Range("K69").Value = 5398.9998
$ws.Range("M69").Value = -4587.9998
$ws.Range("H72").Value = 1974.75
$ws.Range("I72").Value = 1799.6666
$ws.Range("K72").Value = 16196.9994
$ws.Range("M72").Value = -12140.9994
$ws.Range("H131").Value = 697
$ws.Range("J131").Value = 744.55554
$ws.Range("L131").Value = 2233.66662
$ws.Range("N131").Value = -12313.66662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3686207
$ws.Range("I70").Value = 4960
$ws.Range("K70").Value = 4960
$ws.Range("M70").Value = -4690
$ws.Range("H73").Value = 3686207
$ws.Range("I73").Value = 4960
$ws.Range("K73").Value = 4960
$ws.Range("M73").Value = -4024
$ws.Range("H132").Value = 23044.703
$ws.Range("I132").Value = 5395.278
$ws.Range("J132").Value = 58343.555
$ws.Range("K132").Value = 16185.834
$ws.Range("L132").Value = 175030.665
$ws.Range("M132").Value = -13655.834
$ws.Range("N132").Value = -180090.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 328278.28
$ws.Range("I132").Value = 403576.6
$ws.Range("J132").Value = 5571.143
$ws.Range("K132").Value = 1210729.8
$ws.Range("L132").Value = 16713.429
$ws.Range("M132").Value = -1208199.8
$ws.Range("N132").Value = -21773.429
$ws.Range("H136").Value = 2234.0952
$ws.Range("I136").Value = 2105.8
$ws.Range("K136").Value = 6317.400000000001
$ws.Range("M136").Value = -3767.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 3738.75
$ws.Range("I20").Value = 2000
$ws.Range("K20").Value = 2000
$ws.Range("M20").Value = -1760
$ws.Range("H53").Value = 9538
$ws.Range("I53").Value = 9076
$ws.Range("J53").Value = 10000
$ws.Range("K53").Value = 9076
$ws.Range("L53").Value = 10000
$ws.Range("M53").Value = -8469
$ws.Range("N53").Value = -11214
$ws.Range("H81").Value = 2258.3333
$ws.Range("I81").Value = 559.8
$ws.Range("J81").Value = 4381.5
$ws.Range("K81").Value = 1119.6
$ws.Range("L81").Value = 8763
$ws.Range("M81").Value = -58.59999999999991
$ws.Range("N81").Value = -10885
$ws.Range("H84").Value = 2258.3333
$ws.Range("I84").Value = 559.8
$ws.Range("J84").Value = 4381.5
$ws.Range("K84").Value = 5598
$ws.Range("L84").Value = 43815
$ws.Range("M84").Value = -294
$ws.Range("N84").Value = -54423
$ws.Range("H126").Value = 2100.2
$ws.Range("I126").Value = 1714.5238
$ws.Range("J126").Value = 4125
$ws.Range("K126").Value = 5143.5714
$ws.Range("L126").Value = 12375
$ws.Range("M126").Value = -2673.5714
$ws.Range("N126").Value = -17315
$ws.Range("H132").Value = 1914.3684
$ws.Range("I132").Value = 989.5833
$ws.Range("J132").Value = 3499.7144
$ws.Range("K132").Value = 2968.7499
$ws.Range("L132").Value = 10499.1432
$ws.Range("M132").Value = -438.7498999999998
$ws.Range("N132").Value = -15559.1432
